# Add an "eta" column (new column F) to the Sheet1 dataset, shifting the
# existing pt_min..diff columns one to the right, and fill it with the
# per-row eta values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column before column F (the old "pt_min" column becomes G,
# and everything to its right shifts right by one). Insert() on an entire
# column carries formulas/styles/relative refs along with it.
$ws.Range("F1").EntireColumn.Insert()

# Header for the new column.
$ws.Range("F1").Value = "eta"

# New per-row eta values (rows 2-14).
$etaValues = @(0.1, 0.3, 0.5, 0.7, 0.9, 1.1, 1.39, 1.7, 1.9, 2.1, 2.3, 2.54, 2.92)

for ($i = 0; $i -lt $etaValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $etaValues[$i]
}

# Match the recorded selection state after the edit.
$ws.Range("F15").Select()
